$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header row: "_old" -> "_FV2210" and "_new" -> "_FV2304"
$oldHeaders = @("Segmentname_old","Segmentgruppe_old","Segment_old","Datenelement_old","Segment ID_old","Code_old","Qualifier_old","Beschreibung_old","Bedingungsausdruck_old","Bedingung_old")
$newHeaders = @("Segmentname_new","Segmentgruppe_new","Segment_new","Datenelement_new","Segment ID_new","Code_new","Qualifier_new","Beschreibung_new","Bedingungsausdruck_new","Bedingung_new")

for ($i = 0; $i -lt $oldHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = ($oldHeaders[$i] -replace "_old", "_FV2210")
}
for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = ($newHeaders[$i] -replace "_new", "_FV2304")
}

# 2. Turn the used range into an Excel Table ("Table1") spanning A1:U72
$tableRange = $ws.Range("A1:U72")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"

# 3. Freeze the header row (split/freeze at row 2, i.e. top row frozen)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
